$d = $word.ActiveDocument

$old = "Sur l’ensemble de l’année, je suis plutôt satisfait de mon séjour au Canada. Cette année m’aura permis de découvrir ce pays (au moins l’est), ces cultures et aussi de rencontrer de nouvelles personnes. Tout d’abord, d’un point de vue scolaire, en choisissant mes cours, j’ai pu approfondir mes connaissances sur des sujets qui m’intéressaient vraiment"
$new = "En conclusion, je suis vraiment satisfait de mon séjour à l’université de Pittsburgh. Il m’aura permis de découvrir le pays, des cultures et aussi de rencontrer de nouvelles personnes. Tout d’abord, d’un point de vue scolaire, en choisissant mes cours, j’ai pu élargir mes connaissances sur des sujets qui m’intéressaient vraiment"

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
